# Generate Report for Handback
#
# A new handback run has completed for e2e\f9e92b1e-807b-46c9-8a17-88b821d6cb7d.md
# (this is in addition to the existing row, whose own generated artifact names /
# timestamps have also moved forward to the latest run:
# e2e\60f4349b-0f73-4053-91b3-b5a7e9dc8dee.md). Refresh the three report sheets
# (Overview, zh-cn, de-de) accordingly.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newGuid1 = "60f4349b-0f73-4053-91b3-b5a7e9dc8dee"
$newGuid2 = "f9e92b1e-807b-46c9-8a17-88b821d6cb7d"

$zhcnHash   = "c17e4c0197a54b8c2a9b53ea44ff7b723837bd3b"
$dedeHash   = "c17e4c0197a54b8c2a9b53ea44ff7b723837bd3b"
$newRowHash = "5a1583a7d6a4078213ead6144ca99524643607a3"

function Set-Hyperlink($ws, $cell, $url, $display) {
    $rng = $ws.Range($cell)
    if ($rng.Hyperlinks.Count -gt 0) {
        $rng.Hyperlinks.Delete()
    }
    $ws.Hyperlinks.Add($rng, $url, "", "", $display) | Out-Null
}

# ---------------------------------------------------------------------------
# Overview sheet: refresh the existing handback entry, then append the new one
# ---------------------------------------------------------------------------
$overview.Range("A2").Value = "$newGuid1.md"
$overview.Range("B2").Value = "e2e\$newGuid1.md"
$overview.Range("G2").Value = "2016-08-19 21:04:42"
Set-Hyperlink $overview "B2" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1e420eb3e9991c04099401a7b6071cad3f432208/e2e/$newGuid1.md" "e2e\$newGuid1.md"

$oLo = $overview.ListObjects.Item(1)
$oLo.ListRows.Add() | Out-Null
$overview.Range("A3").Value = "$newGuid2.md"
$overview.Range("B3").Value = "e2e\$newGuid2.md"
$overview.Range("C3").Value = ".md"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"
$overview.Range("G3").Value = "2016-08-19 21:04:42"
Set-Hyperlink $overview "B3" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1e420eb3e9991c04099401a7b6071cad3f432208/e2e/$newGuid2.md" "e2e\$newGuid2.md"

# ---------------------------------------------------------------------------
# zh-cn sheet: refresh the existing row's xlf/timestamps, then append new row
# ---------------------------------------------------------------------------
$zhcn.Range("A2").Value = "$newGuid1.md"
$zhcn.Range("G2").Value = "$newGuid1.$zhcnHash.zh-cn.xlf"
$zhcn.Range("H2").Value = "2016-08-19 21:04:37"
$zhcn.Range("I2").Value = "$newGuid1.md"
$zhcn.Range("J2").Value = "$newGuid1.$zhcnHash.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-19 21:04:54"
Set-Hyperlink $zhcn "A2" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1e420eb3e9991c04099401a7b6071cad3f432208/e2e/$newGuid1.md" "$newGuid1.md"
Set-Hyperlink $zhcn "I2" "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/549bd99f278106279ea40ec2a0b69c08d27c016e/e2e/$newGuid1.md" "$newGuid1.md"

$zLo = $zhcn.ListObjects.Item(1)
$zLo.ListRows.Add() | Out-Null
$zhcn.Range("A3").Value = "$newGuid2.md"
$zhcn.Range("B3").Value = ".md"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("D3").Value = "e2e"
$zhcn.Range("E3").Value = "ht"
$zhcn.Range("F3").Value = "True"
$zhcn.Range("G3").Value = "$newGuid2.$newRowHash.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-19 21:04:37"
$zhcn.Range("I3").Value = "$newGuid2.md"
$zhcn.Range("J3").Value = "$newGuid2.$newRowHash.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-19 21:04:54"
$zhcn.Range("L3").Value = ""
$zhcn.Range("M3").Value = "True"
$zhcn.Range("N3").Value = ""
$zhcn.Range("O3").Value = "False"
$zhcn.Range("P3").Value = ""
Set-Hyperlink $zhcn "A3" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1e420eb3e9991c04099401a7b6071cad3f432208/e2e/$newGuid2.md" "$newGuid2.md"
Set-Hyperlink $zhcn "I3" "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/549bd99f278106279ea40ec2a0b69c08d27c016e/e2e/$newGuid2.md" "$newGuid2.md"

# ---------------------------------------------------------------------------
# de-de sheet: refresh the existing row's xlf/timestamps, then append new row
# ---------------------------------------------------------------------------
$dede.Range("A2").Value = "$newGuid1.md"
$dede.Range("G2").Value = "$newGuid1.$dedeHash.de-de.xlf"
$dede.Range("H2").Value = "2016-08-19 21:04:42"
$dede.Range("I2").Value = "$newGuid1.md"
$dede.Range("J2").Value = "$newGuid1.$dedeHash.de-de.xlf"
$dede.Range("K2").Value = "2016-08-19 21:05:02"
Set-Hyperlink $dede "A2" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1e420eb3e9991c04099401a7b6071cad3f432208/e2e/$newGuid1.md" "$newGuid1.md"
Set-Hyperlink $dede "I2" "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/730bc6545aa2409e42964015d5891c0cd52f43b3/e2e/$newGuid1.md" "$newGuid1.md"

$dLo = $dede.ListObjects.Item(1)
$dLo.ListRows.Add() | Out-Null
$dede.Range("A3").Value = "$newGuid2.md"
$dede.Range("B3").Value = ".md"
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("D3").Value = "e2e"
$dede.Range("E3").Value = "ht"
$dede.Range("F3").Value = "True"
$dede.Range("G3").Value = "$newGuid2.$newRowHash.de-de.xlf"
$dede.Range("H3").Value = "2016-08-19 21:04:42"
$dede.Range("I3").Value = "$newGuid2.md"
$dede.Range("J3").Value = "$newGuid2.$newRowHash.de-de.xlf"
$dede.Range("K3").Value = "2016-08-19 21:05:02"
$dede.Range("L3").Value = ""
$dede.Range("M3").Value = "True"
$dede.Range("N3").Value = ""
$dede.Range("O3").Value = "False"
$dede.Range("P3").Value = ""
Set-Hyperlink $dede "A3" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1e420eb3e9991c04099401a7b6071cad3f432208/e2e/$newGuid2.md" "$newGuid2.md"
Set-Hyperlink $dede "I3" "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/730bc6545aa2409e42964015d5891c0cd52f43b3/e2e/$newGuid2.md" "$newGuid2.md"

# ---------------------------------------------------------------------------
# Number formats for date/time columns on the new row (match existing columns)
# ---------------------------------------------------------------------------
$overview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zhcn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zhcn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$dede.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$dede.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
